# Update simulation-derived probability matrix values in-place.
# (added more games, sped up simulate game logic, and drafted optimization logic)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 2).Value = 0.2048611111111111
$ws.Cells.Item(2, 3).Value = 0.5173611111111112
$ws.Cells.Item(2, 10).Value = 0.01736111111111111
$ws.Cells.Item(2, 16).Value = 0.1319444444444444
$ws.Cells.Item(2, 19).Value = 0.1284722222222222

$ws.Cells.Item(3, 2).Value = 0.0130718954248366
$ws.Cells.Item(3, 3).Value = 0.0261437908496732
$ws.Cells.Item(3, 10).Value = 0.006535947712418301
$ws.Cells.Item(3, 16).Value = 0.7124183006535948
$ws.Cells.Item(3, 19).Value = 0.2418300653594771

$ws.Cells.Item(4, 16).Value = 0.7948717948717948
$ws.Cells.Item(4, 19).Value = 0.2051282051282051

$ws.Cells.Item(6, 2).Value = 0.05208333333333334
$ws.Cells.Item(6, 4).Value = 0.01041666666666667
$ws.Cells.Item(6, 6).Value = 0.02604166666666667
$ws.Cells.Item(6, 10).Value = 0.25
$ws.Cells.Item(6, 15).Value = 0.02604166666666667
$ws.Cells.Item(6, 17).Value = 0.2135416666666667
$ws.Cells.Item(6, 18).Value = 0.08333333333333333
$ws.Cells.Item(6, 19).Value = 0.3385416666666667

$ws.Cells.Item(7, 2).Value = 0.1258278145695364
$ws.Cells.Item(7, 4).Value = 0.01986754966887417
$ws.Cells.Item(7, 5).Value = 0.006622516556291391
$ws.Cells.Item(7, 6).Value = 0.04635761589403974
$ws.Cells.Item(7, 10).Value = 0.1258278145695364
$ws.Cells.Item(7, 15).Value = 0.02649006622516556
$ws.Cells.Item(7, 17).Value = 0.152317880794702
$ws.Cells.Item(7, 18).Value = 0.08609271523178808
$ws.Cells.Item(7, 19).Value = 0.4105960264900662

$ws.Cells.Item(8, 2).Value = 0.09411764705882353
$ws.Cells.Item(8, 4).Value = 0.01176470588235294
$ws.Cells.Item(8, 5).Value = 0.004705882352941176
$ws.Cells.Item(8, 6).Value = 0.04470588235294118
$ws.Cells.Item(8, 10).Value = 0.1152941176470588
$ws.Cells.Item(8, 15).Value = 0.01411764705882353
$ws.Cells.Item(8, 17).Value = 0.1811764705882353
$ws.Cells.Item(8, 18).Value = 0.08470588235294117
$ws.Cells.Item(8, 19).Value = 0.4494117647058823

$ws.Cells.Item(9, 2).Value = 0.08520179372197309
$ws.Cells.Item(9, 4).Value = 0.0179372197309417
$ws.Cells.Item(9, 5).Value = 0.004484304932735426
$ws.Cells.Item(9, 6).Value = 0.06278026905829596
$ws.Cells.Item(9, 10).Value = 0.09417040358744394
$ws.Cells.Item(9, 15).Value = 0.03139013452914798
$ws.Cells.Item(9, 17).Value = 0.2331838565022422
$ws.Cells.Item(9, 18).Value = 0.04932735426008968
$ws.Cells.Item(9, 19).Value = 0.42152466367713

$ws.Cells.Item(10, 2).Value = 0.1052202283849918
$ws.Cells.Item(10, 4).Value = 0.02120717781402937
$ws.Cells.Item(10, 5).Value = 0.0008156606851549756
$ws.Cells.Item(10, 6).Value = 0.06525285481239804
$ws.Cells.Item(10, 10).Value = 0.1190864600326264
$ws.Cells.Item(10, 15).Value = 0.01631321370309951
$ws.Cells.Item(10, 17).Value = 0.2398042414355628
$ws.Cells.Item(10, 18).Value = 0.08564437194127243
$ws.Cells.Item(10, 19).Value = 0.3466557911908646

$ws.Cells.Item(11, 7).Value = 0.1440677966101695
$ws.Cells.Item(11, 10).Value = 0.09745762711864407
$ws.Cells.Item(11, 11).Value = 0.2288135593220339
$ws.Cells.Item(11, 12).Value = 0.5169491525423728
$ws.Cells.Item(11, 19).Value = 0.01271186440677966

$ws.Cells.Item(12, 7).Value = 0.7741935483870968
$ws.Cells.Item(12, 10).Value = 0.1774193548387097
$ws.Cells.Item(12, 12).Value = 0.02419354838709677
$ws.Cells.Item(12, 19).Value = 0.02419354838709677

$ws.Cells.Item(13, 7).Value = 0.6764705882352942
$ws.Cells.Item(13, 10).Value = 0.2941176470588235
$ws.Cells.Item(13, 19).Value = 0.02941176470588235

$ws.Cells.Item(14, 7).Value = 0.6
$ws.Cells.Item(14, 19).Value = 0.4

$ws.Cells.Item(15, 6).Value = 0.01673640167364017
$ws.Cells.Item(15, 8).Value = 0.1548117154811715
$ws.Cells.Item(15, 9).Value = 0.07949790794979079
$ws.Cells.Item(15, 10).Value = 0.3807531380753138
$ws.Cells.Item(15, 11).Value = 0.04184100418410042
$ws.Cells.Item(15, 13).Value = 0.01255230125523013
$ws.Cells.Item(15, 15).Value = 0.08368200836820083
$ws.Cells.Item(15, 19).Value = 0.2301255230125523

$ws.Cells.Item(16, 6).Value = 0.01685393258426966
$ws.Cells.Item(16, 8).Value = 0.1348314606741573
$ws.Cells.Item(16, 9).Value = 0.1348314606741573
$ws.Cells.Item(16, 10).Value = 0.4044943820224719
$ws.Cells.Item(16, 11).Value = 0.1067415730337079
$ws.Cells.Item(16, 13).Value = 0.02247191011235955
$ws.Cells.Item(16, 15).Value = 0.0898876404494382
$ws.Cells.Item(16, 19).Value = 0.0898876404494382

$ws.Cells.Item(17, 6).Value = 0.00819672131147541
$ws.Cells.Item(17, 8).Value = 0.1741803278688525
$ws.Cells.Item(17, 9).Value = 0.1045081967213115
$ws.Cells.Item(17, 10).Value = 0.4508196721311475
$ws.Cells.Item(17, 11).Value = 0.0860655737704918
$ws.Cells.Item(17, 13).Value = 0.00819672131147541
$ws.Cells.Item(17, 14).Value = 0.006147540983606557
$ws.Cells.Item(17, 15).Value = 0.05327868852459016
$ws.Cells.Item(17, 19).Value = 0.1086065573770492

$ws.Cells.Item(18, 6).Value = 0.02777777777777778
$ws.Cells.Item(18, 8).Value = 0.2
$ws.Cells.Item(18, 9).Value = 0.09444444444444444
$ws.Cells.Item(18, 10).Value = 0.4388888888888889
$ws.Cells.Item(18, 11).Value = 0.05
$ws.Cells.Item(18, 13).Value = 0.03333333333333333
$ws.Cells.Item(18, 15).Value = 0.06111111111111111
$ws.Cells.Item(18, 19).Value = 0.09444444444444444

$ws.Cells.Item(19, 6).Value = 0.0198961937716263
$ws.Cells.Item(19, 8).Value = 0.2136678200692042
$ws.Cells.Item(19, 9).Value = 0.09688581314878893
$ws.Cells.Item(19, 10).Value = 0.3771626297577855
$ws.Cells.Item(19, 11).Value = 0.09083044982698962
$ws.Cells.Item(19, 13).Value = 0.01557093425605536
$ws.Cells.Item(19, 14).Value = 0.00259515570934256
$ws.Cells.Item(19, 15).Value = 0.08564013840830449
$ws.Cells.Item(19, 19).Value = 0.09775086505190311
